$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated save_data: column G (K = strikeouts) values updated to use
# the actual strikeout count (K) instead of the previous Strike# metric.
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 5
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 6
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 5
$ws.Range("G8").Value = 4
$ws.Range("G9").Value = 3
$ws.Range("G10").Value = 1
$ws.Range("G11").Value = 3
